$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert the new "overview diagram" block right after the
#    "Foreword" M2Doc field paragraph (the m:se.description field),
#    and before the "m:for package | rootPkg.ownedDataPkgs" loop.
# ------------------------------------------------------------------

$targetFieldIndex = 0
for ($i = 1; $i -le $d.Fields.Count; $i++) {
    $code = $d.Fields.Item($i).Code.Text
    if ($code -match "se\.description") {
        $targetFieldIndex = $i
        break
    }
}

$anchorRange = $d.Fields.Item($targetFieldIndex).Code.Paragraphs(1).Range
$anchorRange.Collapse(0)

$newBlockXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:u w:val="single"/>
    </w:rPr>
    <w:br w:type="page"/>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t xml:space="preserve">[CDB] Capella Light </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>Metamodel</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>:</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:fldChar w:fldCharType="begin"/>
  </w:r>
  <w:r>
    <w:instrText xml:space="preserve"> m:'</w:instrText>
  </w:r>
  <w:r>
    <w:instrText>[CDB] Capella Light Metamodel</w:instrText>
  </w:r>
  <w:r>
    <w:instrText xml:space="preserve">'.asImageByRepresentationName().fit(500,650) </w:instrText>
  </w:r>
  <w:r>
    <w:fldChar w:fldCharType="end"/>
  </w:r>
</w:p>
'@

$anchorRange.InsertXML($newBlockXml)

# ------------------------------------------------------------------
# 2) Update the cached PAGE / NUMPAGES field results in the footer
#    (these are plain cached w:t values, not recomputed live here).
# ------------------------------------------------------------------

$footer = $d.Sections(1).Footers(1)
$ffields = $footer.Range.Fields
for ($i = 1; $i -le $ffields.Count; $i++) {
    $fld = $ffields.Item($i)
    $code = $fld.Code.Text
    if ($code -match "NUMPAGES") {
        $fld.Result.Text = "8"
    } elseif ($code -match "PAGE") {
        $fld.Result.Text = "3"
    }
}
